# SECTOR_holdings.xlsx update
# - Bumps the "as of" date in the confidential disclosure note (A9) from
#   2021-03-29 to 2021-03-30.
# - Refreshes the Weight (D) / Percent Change (E) figures for rows 2-6.
#
# The worksheet ships with sheet protection enabled, so locked cells must be
# unprotected before they can be written to, and protection is restored
# afterwards to leave the sheet in a protected state again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Disclosure footer text (shared string used by cell A9) - just the date changes.
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-30 for illustrative purposes only and are subject to change."

# Row 2 (XLB / Materials Select Sector SPDR)
$ws.Range("D2").Value = 0.2546514430736706
$ws.Range("E2").Value = -0.002895995970788268

# Row 3 (XLF / Financial Select Sector SPDR Fund)
$ws.Range("D3").Value = 0.2490786600865685
$ws.Range("E3").Value = 0.007044320516583547

# Row 4 (XLK / Technology Select Sector SPDR Fund)
$ws.Range("D4").Value = 0.2509149965881966
$ws.Range("E4").Value = -0.009471131989695403

# Row 5 (XLC / Communication Services Select Sector SPDR Fund)
$ws.Range("D5").Value = 0.2453549002515643
$ws.Range("E5").Value = 0

# Row 6 (Total)
$ws.Range("E6").Value = -0.00135932868848665

# Restore protection on the sheet.
$ws.Protect()
